$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.831.55'
$ws.Range('E2').Value = '  +5.26%  '
$ws.Range('D3').Value = '3.289.75'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '627.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.417'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.713'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.33%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = '3.287.22'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.598'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000275'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.180'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('D16').Value = '3.877.32'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '90.429.27'
$ws.Range('E17').Value = '  +5.09%  '
$ws.Range('D18').Value = '3.272.60'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000190'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +48.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').Value = '3.428.76'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '76.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.177'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.994'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '566.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.26'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.50%  '
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +24.35%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.135'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.395'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '183.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.90%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '149.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '44.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.130'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.00%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.630'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.36%  '
